$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2
$ws.Range("D2").Value = '27.214.32'

# Row 3: update D3, E3
$ws.Range("D3").Value = '1.685.88'
$ws.Range("E3").Value = '  +0.57%  '

# Row 4: update E4
$ws.Range("E4").Value = '  -0.10%  '

# Row 5: update D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.24%  '

# Row 6: update D6, E6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.519'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.29%  '

# Row 7: update E7
$ws.Range("E7").Value = '  -0.11%  '

# Row 9: update D9, E9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.80'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.28%  '

# Row 10: update D10, E10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0625'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.85%  '

# Row 11: update E11
$ws.Range("E11").Value = '  +0.17%  '

# Row 12: update D12, E12
$ws.Range("D12").Value = '1.922.92'
$ws.Range("E12").Value = '  +0.53%  '

# Row 13: update D13, E13
$ws.Range("D13").Value = '1.684.41'
$ws.Range("E13").Value = '  +0.65%  '

# Row 14: update D14, E14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.97%  '

# Row 15: update D15, E15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.547'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.49%  '

# Row 16: update D16, E16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.02%  '

# Row 17: update D17, E17
$ws.Range("D17").Value = '27.208.37'
$ws.Range("E17").Value = '  +0.80%  '

# Row 18: update D18, E18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '239.67'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.18%  '

# Row 19: update D19, E19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.80%  '

# Row 20: update D20, E20
$ws.Range("D20").Value = '0.0₃0744'
$ws.Range("E20").Value = '  +1.49%  '

# Row 21: update E21
$ws.Range("E21").Value = '  -0.12%  '

# Row 22: update D22, E22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.57'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.85%  '

# Row 23: update D23, E23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.64%  '

# Row 24: update E24
$ws.Range("E24").Value = '  -3.15%  '

# Row 25: update D25, E25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.34'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.84%  '

# Row 26: update D26, E26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.29'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.84%  '

# Row 27: update D27, E27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.39'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.97%  '

# Row 28: update E28
$ws.Range("E28").Value = '  +0.87%  '

# Row 29: update D29, E29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.17%  '

# Row 30: update D30, E30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0500'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.34%  '

# Row 31: update E31
$ws.Range("E31").Value = '  +0.45%  '

# Row 32: update D32, E32
$ws.Range("D32").Value = '1.578.73'
$ws.Range("E32").Value = '  +6.27%  '

# Row 33: update E33
$ws.Range("E33").Value = '  +1.70%  '

# Row 34: update E34
$ws.Range("E34").Value = '  +2.59%  '

# Row 35: update E35
$ws.Range("E35").Value = '  +0.43%  '

# Row 36: update B36, C36, D36, E36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.954'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.84%  '

# Row 37: update B37, C37, D37, E37
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.602'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.99%  '

# Row 38: update E38
$ws.Range("E38").Value = '  -1.01%  '

# Row 39: update E39
$ws.Range("E39").Value = '  -0.58%  '

# Row 40: update D40, E40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.00%  '

# Row 41: update D41, E41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.40'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.82%  '

# Row 42: update E42
$ws.Range("E42").Value = '  -0.10%  '

# Row 43: update E43
$ws.Range("E43").Value = '  -4.00%  '

# Row 44: update E44
$ws.Range("E44").Value = '  -2.47%  '

# Row 45: update D45, E45
$ws.Range("D45").Value = '1.831.97'
$ws.Range("E45").Value = '  +0.73%  '

# Row 46: update D46, E46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.787'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.50%  '

# Row 47: update D47, E47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.99'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.39%  '

# Row 48: update E48
$ws.Range("E48").Value = '  +4.11%  '

# Row 49: update D49, E49
$ws.Range("D49").Value = '0.0₆0106'
$ws.Range("E49").Value = '  -0.13%  '

# Row 50: update B50, C50, D50, E50
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.104'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.15%  '

# Row 51: update B51, C51, D51, E51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.71%  '
